# Generate Report for Handoff
#
# Updates the "localization-status" workbook to reflect a freshly
# generated handback/handoff report for the
# d7f87c26-6c67-49cd-9b65-d752389df486 item (last row, row 7, of each
# sheet):
#   - Overview!G7              (Latest HO Xliff Generate Date)
#   - zh-cn!H7                 (Latest Handback DateTime)
#   - de-de!H7                 (Latest Handback DateTime)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-20 14:45:30"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-20 14:45:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-20 14:45:30"
